$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.593300805548547
$ws.Range("C2").Value = 0.798316168795588
$ws.Range("D2").Value = 0.809737487250624
$ws.Range("E2").Value = 0.801269845163997
$ws.Range("F2").Value = 0.802759474182296
$ws.Range("G2").Value = 0.803503442074255

$ws.Range("B3").Value = 0.891425446185947
$ws.Range("C3").Value = 0.734343884976076
$ws.Range("D3").Value = 0.66011845184374
$ws.Range("E3").Value = 0.773993520808291
$ws.Range("F3").Value = 0.769785232844933
$ws.Range("G3").Value = 0.734188607563879

$ws.Range("B4").Value = 0.837842434332639
$ws.Range("C4").Value = 0.649556135900333
$ws.Range("D4").Value = 0.544523480585528
$ws.Range("E4").Value = 0.727832934920089
$ws.Range("F4").Value = 0.722448306363765
$ws.Range("G4").Value = 0.676610115099755

$ws.Range("B5").Value = 0.747059011250423
$ws.Range("C5").Value = 0.760748646394783
$ws.Range("D5").Value = 0.820640960368556
$ws.Range("E5").Value = 0.745104631264865
$ws.Range("F5").Value = 0.72988960058406
$ws.Range("G5").Value = 0.728092737150705

$ws.Range("B6").Value = 0.620917729383005
$ws.Range("C6").Value = 0.78913119672879
$ws.Range("D6").Value = 0.855042207071307
$ws.Range("E6").Value = 0.753586689991038
$ws.Range("F6").Value = 0.834664561377067
$ws.Range("G6").Value = 0.754698076525027

$ws.Range("B7").Value = 0.701622795907415
$ws.Range("C7").Value = 0.712710646471761
$ws.Range("D7").Value = 0.757480297774307
$ws.Range("E7").Value = 0.637897645695564
$ws.Range("F7").Value = 0.877803125109165
$ws.Range("G7").Value = 0.700890373546844

$ws.Range("B8").Value = 0.71036942310177
$ws.Range("C8").Value = 0.767459227638374
$ws.Range("D8").Value = 0.769430442315501
$ws.Range("E8").Value = 0.758346955137085
$ws.Range("F8").Value = 0.84096493071521
$ws.Range("G8").Value = 0.746768125063093

$ws.Range("B9").Value = 0.776451040625449
$ws.Range("C9").Value = 0.752356847825706
$ws.Range("D9").Value = 0.829686492804909
$ws.Range("E9").Value = 0.705157484400394
$ws.Range("F9").Value = 0.849273919987832
$ws.Range("G9").Value = 0.66342378640839

$ws.Range("B10").Value = 0.869708886859866
$ws.Range("C10").Value = 0.815652162844027
$ws.Range("D10").Value = 0.834311730887263
$ws.Range("E10").Value = 0.8156709550428
$ws.Range("F10").Value = 0.867891452365697
$ws.Range("G10").Value = 0.7671375817381

$ws.Range("B11").Value = 0.849616378214432
$ws.Range("C11").Value = 0.734028088236682
$ws.Range("D11").Value = 0.756084447572338
$ws.Range("E11").Value = 0.685545279841557
$ws.Range("F11").Value = 0.842706846604206
$ws.Range("G11").Value = 0.69712559445501

$ws.Range("B12").Value = 0.891121217137727
$ws.Range("C12").Value = 0.819983433106146
$ws.Range("D12").Value = 0.87483983315479
$ws.Range("E12").Value = 0.673226727255466
$ws.Range("F12").Value = 0.913082349886563
$ws.Range("G12").Value = 0.869575896966582
